$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 19; existing rows 19:52 shift down to 21:54
$ws.Range("A19:A20").EntireRow.Insert()

# New row 19 data
$ws.Cells.Item(19, 1).Value = 1
$ws.Cells.Item(19, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(19, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(19, 4).Value = 44607
$ws.Cells.Item(19, 4).Style = $ws.Cells.Item(21, 4).Style
$ws.Cells.Item(19, 4).NumberFormat = $ws.Cells.Item(21, 4).NumberFormat
$ws.Cells.Item(19, 5).Value = 15
$ws.Cells.Item(19, 6).Value = 100112009
$ws.Cells.Item(19, 7).Value = "Acelga"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 160
$ws.Cells.Item(19, 11).Value = 1800
$ws.Cells.Item(19, 12).Value = 2000
$ws.Cells.Item(19, 13).Value = 1900
$ws.Cells.Item(19, 14).Value = "`$/atado 2,5 a 3 kilos"
$ws.Cells.Item(19, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(19, 16).Value = 633
$ws.Cells.Item(19, 17).Value = 3
$ws.Cells.Item(19, 18).Value = "Hortaliza"

# New row 20 data
$ws.Cells.Item(20, 1).Value = 1
$ws.Cells.Item(20, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(20, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(20, 4).Value = 44607
$ws.Cells.Item(20, 4).Style = $ws.Cells.Item(21, 4).Style
$ws.Cells.Item(20, 4).NumberFormat = $ws.Cells.Item(21, 4).NumberFormat
$ws.Cells.Item(20, 5).Value = 15
$ws.Cells.Item(20, 6).Value = 100112009
$ws.Cells.Item(20, 7).Value = "Acelga"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Segunda"
$ws.Cells.Item(20, 10).Value = 250
$ws.Cells.Item(20, 11).Value = 1300
$ws.Cells.Item(20, 12).Value = 1500
$ws.Cells.Item(20, 13).Value = 1400
$ws.Cells.Item(20, 14).Value = "`$/atado 2,5 a 3 kilos"
$ws.Cells.Item(20, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(20, 16).Value = 467
$ws.Cells.Item(20, 17).Value = 3
$ws.Cells.Item(20, 18).Value = "Hortaliza"
